$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "CLS-8929"
$ws.Range("B2").Value = 2.9
$ws.Range("C2").Value = "ANO"

$ws.Range("A3").Select()
